# Insert a new data row at row 110 (this shifts the existing rows 110..161
# down to 111..162, matching the diff which shows every row from 110
# onward moving down by one and a brand-new row appearing at the top of
# that range with fresh data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new record's values.
$ws.Cells.Item(110, 1).Value2  = 11
$ws.Cells.Item(110, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(110, 3).Value2  = "Bíobío"
$ws.Cells.Item(110, 4).Value2  = 44839
$ws.Cells.Item(110, 5).Value2  = 8
$ws.Cells.Item(110, 6).Value2  = "Fruta"
$ws.Cells.Item(110, 7).Value2  = 100102
$ws.Cells.Item(110, 8).Value2  = "Cítricos"
$ws.Cells.Item(110, 9).Value2  = 100102004
$ws.Cells.Item(110, 10).Value2 = "Mandarina"
$ws.Cells.Item(110, 11).Value2 = "Murcott"
$ws.Cells.Item(110, 12).Value2 = "Primera"
$ws.Cells.Item(110, 13).Value2 = 300
$ws.Cells.Item(110, 14).Value2 = 6500
$ws.Cells.Item(110, 15).Value2 = 7000
$ws.Cells.Item(110, 16).Value2 = 6750
$ws.Cells.Item(110, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(110, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(110, 19).Value2 = 375
$ws.Cells.Item(110, 20).Value2 = 18
